# requirements.txt addition - G15
#
# 1) Convert F101:F109 from text ("inlineStr") to real numbers (same values).
# 2) Append 18 new data rows (110-127), duplicating the last two "scrape"
#    batches (timestamps 13:15:12 and 13:20:26) for Disney's Animal Kingdom.
#    Rows 110-118 store column F as numbers; rows 119-127 store column F as text.
#    Columns A-E always stay plain text (dates/times must not be auto-converted
#    into Excel date/number serials).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-NumCell($row, $col, $value) {
    $ws.Cells.Item($row, $col).Value = $value
}

function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# ---------------------------------------------------------------------------
# Step 1: fix F101:F109 - change the existing text values into real numbers
# ---------------------------------------------------------------------------
$existingF = @{
    101 = 105
    102 = 50
    103 = 40
    104 = 15
    105 = 50
    106 = 25
    107 = 35
    108 = 80
    109 = 5
}
foreach ($r in $existingF.Keys) {
    Set-NumCell $r 6 $existingF[$r]
}

# ---------------------------------------------------------------------------
# Step 2: append new rows 110-127
# ---------------------------------------------------------------------------
# Each entry: row, date, time, park, area, name, wait(F)
$newRows = @(
    @(110, "2024-07-25", "13:15:12", "Disney's Animal Kingdom", "Pandora - The World of Avatar", "Avatar Flight of Passage", 105),
    @(111, "2024-07-25", "13:15:12", "Disney's Animal Kingdom", "Dinoland USA", "DINOSAUR", 45),
    @(112, "2024-07-25", "13:15:12", "Disney's Animal Kingdom", "Asia", "Expedition Everest", 40),
    @(113, "2024-07-25", "13:15:12", "Disney's Animal Kingdom", "Discovery Island", "It's Tough to be a Bug!", 10),
    @(114, "2024-07-25", "13:15:12", "Disney's Animal Kingdom", "Asia", "Kali River Rapids", 50),
    @(115, "2024-07-25", "13:15:12", "Disney's Animal Kingdom", "Africa", "Kilimanjaro Safaris", 25),
    @(116, "2024-07-25", "13:15:12", "Disney's Animal Kingdom", "Discovery Island", "Meet Favorite Disney Pals at Adventurers Outpost", 15),
    @(117, "2024-07-25", "13:15:12", "Disney's Animal Kingdom", "Pandora - The World of Avatar", "Na'vi River Journey", 85),
    @(118, "2024-07-25", "13:15:12", "Disney's Animal Kingdom", "Dinoland USA", "TriceraTop Spin", 5)
)

$newRowsText = @(
    @(119, "2024-07-25", "13:20:26", "Disney's Animal Kingdom", "Pandora - The World of Avatar", "Avatar Flight of Passage", "105"),
    @(120, "2024-07-25", "13:20:26", "Disney's Animal Kingdom", "Dinoland USA", "DINOSAUR", "45"),
    @(121, "2024-07-25", "13:20:26", "Disney's Animal Kingdom", "Asia", "Expedition Everest", "40"),
    @(122, "2024-07-25", "13:20:26", "Disney's Animal Kingdom", "Discovery Island", "It's Tough to be a Bug!", "10"),
    @(123, "2024-07-25", "13:20:26", "Disney's Animal Kingdom", "Asia", "Kali River Rapids", "50"),
    @(124, "2024-07-25", "13:20:26", "Disney's Animal Kingdom", "Africa", "Kilimanjaro Safaris", "35"),
    @(125, "2024-07-25", "13:20:26", "Disney's Animal Kingdom", "Discovery Island", "Meet Favorite Disney Pals at Adventurers Outpost", "15"),
    @(126, "2024-07-25", "13:20:26", "Disney's Animal Kingdom", "Pandora - The World of Avatar", "Na'vi River Journey", "70"),
    @(127, "2024-07-25", "13:20:26", "Disney's Animal Kingdom", "Dinoland USA", "TriceraTop Spin", "5")
)

foreach ($entry in $newRows) {
    $r = $entry[0]
    Set-TextCell $r 1 $entry[1]
    Set-TextCell $r 2 $entry[2]
    Set-TextCell $r 3 $entry[3]
    Set-TextCell $r 4 $entry[4]
    Set-TextCell $r 5 $entry[5]
    Set-NumCell $r 6 $entry[6]
}

foreach ($entry in $newRowsText) {
    $r = $entry[0]
    Set-TextCell $r 1 $entry[1]
    Set-TextCell $r 2 $entry[2]
    Set-TextCell $r 3 $entry[3]
    Set-TextCell $r 4 $entry[4]
    Set-TextCell $r 5 $entry[5]
    Set-TextCell $r 6 $entry[6]
}
